$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$new = $wb.Worksheets.Add()
$wb.Worksheets.Item("Sheet1").Delete()

$ws = $wb.Worksheets.Item(1)
$ws.Range("A1").Value = "Code"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Region"
$ws.Range("D1").Value = "MaxWorkingHourMonthly"

$header = $ws.Range("A1:C1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4108
$header.NumberFormat = "@"

$d1 = $ws.Range("D1")
$d1.Font.Bold = $true
$d1.HorizontalAlignment = -4108
$d1.VerticalAlignment = -4108

$ws.Columns.Item(2).ColumnWidth = 31.42578125
$ws.Columns.Item(4).ColumnWidth = 24.28515625

Write-Host "Done"
foreach ($s in $wb.Worksheets) { Write-Host $s.Name }
